# Update cryptocurrency price/volume data (and the Filecoin/OKB row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'90.305.00"
$ws.Range('E2').Value = "'  +0.31%  "
$ws.Range('D3').Value = "'3.163.18"
$ws.Range('E3').Value = "'  +2.59%  "
$ws.Range('E4').Value = "'  +0.04%  "
$ws.Range('D5').Value = "'237.71"
$ws.Range('E5').Value = "'  +2.11%  "
$ws.Range('D6').Value = "'621.75"
$ws.Range('E6').Value = "'  +0.66%  "
$ws.Range('D7').Value = "'1.10"
$ws.Range('E7').Value = "'  +4.91%  "
$ws.Range('E8').Value = "'  +2.99%  "
$ws.Range('E9').Value = "'  -0.02%  "
$ws.Range('D10').Value = "'3.163.63"
$ws.Range('E10').Value = "'  +2.66%  "
$ws.Range('D11').Value = "'0.739"
$ws.Range('E11').Value = "'  +2.14%  "
$ws.Range('E12').Value = "'  +3.78%  "
$ws.Range('D13').Value = "'0.0000246"
$ws.Range('E13').Value = "'  +0.31%  "
$ws.Range('D14').Value = "'35.34"
$ws.Range('E14').Value = "'  +0.65%  "
$ws.Range('D15').Value = "'5.53"
$ws.Range('E15').Value = "'  +2.51%  "
$ws.Range('D16').Value = "'90.474.87"
$ws.Range('E16').Value = "'  +0.60%  "
$ws.Range('D17').Value = "'3.764.91"
$ws.Range('D18').Value = "'3.183.01"
$ws.Range('E18').Value = "'  +3.34%  "
$ws.Range('D19').Value = "'3.69"
$ws.Range('E19').Value = "'  -4.29%  "
$ws.Range('D20').Value = "'15.05"
$ws.Range('E20').Value = "'  +8.78%  "
$ws.Range('D21').Value = "'5.86"
$ws.Range('E21').Value = "'  +7.18%  "
$ws.Range('D22').Value = "'0.0000204"
$ws.Range('E22').Value = "'  -4.04%  "
$ws.Range('D23').Value = "'440.54"
$ws.Range('E23').Value = "'  +2.05%  "
$ws.Range('D24').Value = "'9.10"
$ws.Range('E24').Value = "'  +3.79%  "
$ws.Range('D25').Value = "'5.74"
$ws.Range('E25').Value = "'  +0.17%  "
$ws.Range('D26').Value = "'89.01"
$ws.Range('E26').Value = "'  +3.03%  "
$ws.Range('D27').Value = "'11.98"
$ws.Range('E27').Value = "'  +1.16%  "
$ws.Range('D28').Value = "'3.329.89"
$ws.Range('E28').Value = "'  +2.89%  "
$ws.Range('D30').Value = "'0.127"
$ws.Range('E30').Value = "'  +46.39%  "
$ws.Range('D31').Value = "'0.230"
$ws.Range('E31').Value = "'  +18.93%  "
$ws.Range('D32').Value = "'0.169"
$ws.Range('E32').Value = "'  +7.70%  "
$ws.Range('D33').Value = "'9.47"
$ws.Range('E33').Value = "'  +3.56%  "
$ws.Range('D34').Value = "'0.998"
$ws.Range('E34').Value = "'  +0.02%  "
$ws.Range('E35').Value = "'  +12.01%  "
$ws.Range('D36').Value = "'7.80"
$ws.Range('E36').Value = "'  +10.60%  "
$ws.Range('D37').Value = "'26.31"
$ws.Range('E37').Value = "'  +2.79%  "
$ws.Range('D38').Value = "'504.43"
$ws.Range('E38').Value = "'  +0.74%  "
$ws.Range('D39').Value = "'1.95"
$ws.Range('E39').Value = "'  +3.89%  "
$ws.Range('D40').Value = "'1.35"
$ws.Range('E40').Value = "'  +7.24%  "
$ws.Range('D41').Value = "'0.450"
$ws.Range('E41').Value = "'  +12.22%  "
$ws.Range('D42').Value = "'3.75"
$ws.Range('E42').Value = "'  +3.78%  "
$ws.Range('D43').Value = "'3.41"
$ws.Range('E43').Value = "'  -9.21%  "
$ws.Range('D44').Value = "'22.10"
$ws.Range('E44').Value = "'  -0.06%  "
$ws.Range('E45').Value = "'  -0.01%  "
$ws.Range('D46').Value = "'0.721"
$ws.Range('E46').Value = "'  +6.77%  "
$ws.Range('E47').Value = "'  +4.08%  "
$ws.Range('E48').Value = "'  +2.87%  "
$ws.Range('D49').Value = "'1.37"
$ws.Range('E49').Value = "'  +4.82%  "
$ws.Range('B50').Value = "'OKB"
$ws.Range('C50').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('D50').Value = "'44.02"
$ws.Range('E50').Value = "'  -0.99%  "
$ws.Range('B51').Value = "'Filecoin"
$ws.Range('C51').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('D51').Value = "'4.41"
$ws.Range('E51').Value = "'  +1.43%  "
